$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new "Area" label + value
$ws.Range("Q3").Value = "Area"
$ws.Range("R3").Value = 257.69726249999991

# Units labels first (matches the order they were typed by the author)
$ws.Range("S9").Value = "Mpa"
$ws.Range("S10").Value = "Gpa"

# Row 9: Exx FEM (numerical) result, in MPa
$ws.Range("Q9").Value = "Exx FEM"
$ws.Range("R9").Formula = "=9000/(R3*R6)"

# Row 10: same result converted to GPa
$ws.Range("R10").Formula = "=R9*10^-3"

# Row 11: Exx Teoria Classica (analytical) result, in MPa
$ws.Range("Q11").Value = "Exx Teoria Clássica"
$ws.Range("R11").Formula = "=9000/(0.0012356*R3)"
$ws.Range("S11").Value = "Mpa"

# Row 12: same result converted to GPa
$ws.Range("R12").Formula = "=R11*10^-3"
$ws.Range("S12").Value = "Gpa"

# Match the final selection / active cell left by the author
$ws.Range("R12").Select() | Out-Null

$wb.Save()
